$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Existing rows 2-10 keep their Key (A) values, only some Value (B) columns and
# the set of rows below them change. Rewrite the whole table to match target.

$data = @(
    @("title", "Project Bloom: A Quest for Home"),
    @("credits_desc", "Made by: RENEGADEWARE"),
    @("options", "OPTIONS"),
    @("music", "MUSIC"),
    @("sound", "SOUND"),
    @("speech", "SPEECH"),
    @("close", "CLOSE"),
    @("on", "ON"),
    @("off", "OFF"),
    @("season_winter", "Winter"),
    @("season_spring", "Spring"),
    @("season_summer", "Summer"),
    @("season_autumn", "Autumn"),
    @("atmosphere_altitude", "Altitude"),
    @("atmosphere_humidity", "Humidity"),
    @("atmosphere_temperature", "Temperature"),
    @("atmosphere_windStrength", "Wind Strength"),
    @("climate_temperate", "Temperate"),
    @("region_NA", "North American Great Plains")
)

# Remove the old "ready"/"READY"/"go"/"GO" rows that used to live at rows 11-12
# by rewriting the full block starting at row 2.
$row = 2
foreach ($pair in $data) {
    $ws.Cells.Item($row, 1).Value = $pair[0]
    $cellB = $ws.Cells.Item($row, 2)
    $cellB.Value = $pair[1]
    $cellB.Style = "Normal"
    $cellB.WrapText = $true
    $row = $row + 1
}

$ws.Range("A20").Select()
